# Update countries & provincias Spain
# - Fix Grecia / Noruega ordering (Grecia now sorts before Noruega) and
#   refresh Grecia's COVID figures for the day; Noruega keeps its previous
#   figures (just shifted down one row).
# - Refresh several other countries' COVID figures.
# - Bump the "Datos actualizados" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Grecia / Noruega (row 91 becomes Grecia w/ new data, row 92 becomes
#     Noruega w/ the data that used to be on row 91) ---
$ws.Cells.Item(91, 1).Value = "Grecia"
$ws.Cells.Item(91, 2).Value = 11386
$ws.Cells.Item(91, 3).Value = 186
$ws.Cells.Item(91, 4).Value = 3804
$ws.Cells.Item(91, 5).Value = 7302
$ws.Cells.Item(91, 7).Value = 1
$ws.Cells.Item(91, 8).Value = 280

$ws.Cells.Item(92, 1).Value = "Noruega"
$ws.Cells.Item(92, 2).Value = 11254
$ws.Cells.Item(92, 3).Value = 23
$ws.Cells.Item(92, 4).Value = 9348
$ws.Cells.Item(92, 5).Value = 1642
$ws.Cells.Item(92, 8).Value = 264

# --- Refresh other countries' figures ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 6405788
$ws.Cells.Item(4, 3).Value = 16731
$ws.Cells.Item(4, 4).Value = 3639838
$ws.Cells.Item(4, 5).Value = 2573570
$ws.Cells.Item(4, 7).Value = 269
$ws.Cells.Item(4, 8).Value = 192380

# Row 5: India
$ws.Cells.Item(5, 2).Value = 4103694
$ws.Cells.Item(5, 3).Value = 83455
$ws.Cells.Item(5, 4).Value = 3172300
$ws.Cells.Item(5, 5).Value = 860761
$ws.Cells.Item(5, 7).Value = 998
$ws.Cells.Item(5, 8).Value = 70633

# Row 6: Brasil
$ws.Cells.Item(6, 2).Value = 4093586
$ws.Cells.Item(6, 3).Value = 1785
$ws.Cells.Item(6, 5).Value = 689684
$ws.Cells.Item(6, 7).Value = 75
$ws.Cells.Item(6, 8).Value = 125659

# Row 21: Turquia
$ws.Cells.Item(21, 2).Value = 278228
$ws.Cells.Item(21, 3).Value = 1673
$ws.Cells.Item(21, 4).Value = 250092
$ws.Cells.Item(21, 5).Value = 21516
$ws.Cells.Item(21, 7).Value = 56
$ws.Cells.Item(21, 8).Value = 6620

# Row 29: Israel
$ws.Cells.Item(29, 2).Value = 128768
$ws.Cells.Item(29, 3).Value = 2349
$ws.Cells.Item(29, 4).Value = 101478
$ws.Cells.Item(29, 5).Value = 26283
$ws.Cells.Item(29, 7).Value = 14
$ws.Cells.Item(29, 8).Value = 1007

# Row 57: Argelia
$ws.Cells.Item(57, 2).Value = 46071
$ws.Cells.Item(57, 3).Value = 298
$ws.Cells.Item(57, 4).Value = 32481
$ws.Cells.Item(57, 5).Value = 12041
$ws.Cells.Item(57, 7).Value = 10
$ws.Cells.Item(57, 8).Value = 1549

# Row 70: Irlanda
$ws.Cells.Item(70, 2).Value = 29534
$ws.Cells.Item(70, 3).Value = 231
$ws.Cells.Item(70, 5).Value = 4393

# Row 72: Chequia
$ws.Cells.Item(72, 2).Value = 27560
$ws.Cells.Item(72, 3).Value = 311
$ws.Cells.Item(72, 4).Value = 19033
$ws.Cells.Item(72, 5).Value = 8096
$ws.Cells.Item(72, 7).Value = 2
$ws.Cells.Item(72, 8).Value = 431

# Row 75: Estado de Palestina
$ws.Cells.Item(75, 5).Value = 8737
$ws.Cells.Item(75, 7).Value = 7
$ws.Cells.Item(75, 8).Value = 177

# Row 130: Sri Lanka
$ws.Cells.Item(130, 2).Value = 3121
$ws.Cells.Item(130, 3).Value = 6
$ws.Cells.Item(130, 5).Value = 191

# Row 138: Bahamas
$ws.Cells.Item(138, 2).Value = 2476
$ws.Cells.Item(138, 3).Value = 90
$ws.Cells.Item(138, 4).Value = 948
$ws.Cells.Item(138, 5).Value = 1473
$ws.Cells.Item(138, 7).Value = 5
$ws.Cells.Item(138, 8).Value = 55

# Row 153: Republica de Chipre
$ws.Cells.Item(153, 2).Value = 1507
$ws.Cells.Item(153, 3).Value = 5
$ws.Cells.Item(153, 5).Value = 249

# --- Bump the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 19:21"
